$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.0609700120985508
$ws.Cells.Item(2, 2).Value = 0.9848241806030273
$ws.Cells.Item(2, 3).Value = 0.01865122653543949
$ws.Cells.Item(2, 4).Value = 0.9973869323730469
$ws.Cells.Item(3, 1).Value = 0.009699323214590549
$ws.Cells.Item(3, 2).Value = 0.998573362827301
$ws.Cells.Item(3, 3).Value = 0.009068938903510571
$ws.Cells.Item(3, 4).Value = 0.9977889657020569
$ws.Cells.Item(4, 1).Value = 0.005020159762352705
$ws.Cells.Item(4, 2).Value = 0.9988051652908325
$ws.Cells.Item(4, 3).Value = 0.004236927721649408
$ws.Cells.Item(4, 4).Value = 0.9981909394264221
$ws.Cells.Item(5, 1).Value = 0.002538583241403103
$ws.Cells.Item(5, 2).Value = 0.9992867112159729
$ws.Cells.Item(5, 3).Value = 0.0008564945892430842
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(6, 1).Value = 0.001806184882298112
$ws.Cells.Item(6, 2).Value = 0.9995898604393005
$ws.Cells.Item(6, 3).Value = 0.0004952636081725359
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(7, 1).Value = 0.00122053746599704
$ws.Cells.Item(7, 2).Value = 0.9997503161430359
$ws.Cells.Item(7, 3).Value = 0.0003024678735528141
$ws.Cells.Item(7, 4).Value = 0.9997990131378174
$ws.Cells.Item(8, 1).Value = 0.001420272863470018
$ws.Cells.Item(8, 2).Value = 0.9996968507766724
$ws.Cells.Item(8, 3).Value = [double]"8.746654202695936E-05"
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(9, 1).Value = 0.0007689269259572029
$ws.Cells.Item(9, 2).Value = 0.9997860193252563
$ws.Cells.Item(9, 3).Value = 0.0001022857177304104
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(10, 1).Value = 0.00162713963072747
$ws.Cells.Item(10, 2).Value = 0.9997146725654602
$ws.Cells.Item(10, 3).Value = [double]"8.391196752199903E-05"
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(11, 1).Value = 0.0006352785276249051
$ws.Cells.Item(11, 2).Value = 0.9998930096626282
$ws.Cells.Item(11, 3).Value = [double]"3.402368020033464E-05"
$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(12, 1).Value = 0.0006416210671886802
$ws.Cells.Item(12, 2).Value = 0.9998394846916199
$ws.Cells.Item(12, 3).Value = [double]"3.27190755342599E-05"
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(13, 1).Value = 0.001049440004862845
$ws.Cells.Item(13, 2).Value = 0.9997681975364685
$ws.Cells.Item(13, 3).Value = [double]"2.287676215928514E-05"
$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(14, 1).Value = 0.0003085793869104236
$ws.Cells.Item(14, 2).Value = 0.999910831451416
$ws.Cells.Item(14, 3).Value = [double]"9.450728248339146E-05"
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(15, 1).Value = 0.0009080973104573786
$ws.Cells.Item(15, 2).Value = 0.9998038411140442
$ws.Cells.Item(15, 3).Value = [double]"5.411836809798842E-06"
$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(16, 1).Value = 0.0004154359921813011
$ws.Cells.Item(16, 2).Value = 0.9998573660850525
$ws.Cells.Item(16, 3).Value = [double]"3.545280151229235E-06"
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(17, 1).Value = 0.0006451503140851855
$ws.Cells.Item(17, 2).Value = 0.9999643564224243
$ws.Cells.Item(17, 3).Value = [double]"5.901569693378406E-06"
$ws.Cells.Item(17, 4).Value = 1
$ws.Cells.Item(18, 1).Value = 0.0004152987676206976
$ws.Cells.Item(18, 2).Value = 0.9998573660850525
$ws.Cells.Item(18, 3).Value = [double]"4.173203706159256E-06"
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(19, 1).Value = 0.000553894555196166
$ws.Cells.Item(19, 2).Value = 0.999910831451416
$ws.Cells.Item(19, 3).Value = [double]"3.872068248256255E-07"
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(20, 1).Value = 0.0004408551612868905
$ws.Cells.Item(20, 2).Value = 0.9998930096626282
$ws.Cells.Item(20, 3).Value = [double]"8.610165309619333E-07"
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(21, 1).Value = 0.0003552399866748601
$ws.Cells.Item(21, 2).Value = 0.9998930096626282
$ws.Cells.Item(21, 3).Value = [double]"3.56089401520876E-07"
$ws.Cells.Item(21, 4).Value = 1
$ws.Cells.Item(22, 1).Value = 0.0003149453259538859
$ws.Cells.Item(22, 2).Value = 0.9999464750289917
$ws.Cells.Item(22, 3).Value = [double]"1.388690975545614E-07"
$ws.Cells.Item(22, 4).Value = 1
$ws.Cells.Item(23, 1).Value = 0.0002008400042541325
$ws.Cells.Item(23, 2).Value = 0.9999821782112122
$ws.Cells.Item(23, 3).Value = [double]"6.256043860730642E-08"
$ws.Cells.Item(23, 4).Value = 1
$ws.Cells.Item(24, 1).Value = 0.0001034078450175002
$ws.Cells.Item(24, 2).Value = 0.9999643564224243
$ws.Cells.Item(24, 3).Value = [double]"6.435783461711253E-08"
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(25, 1).Value = 0.0001139565647463314
$ws.Cells.Item(25, 2).Value = 0.9999643564224243
$ws.Cells.Item(25, 3).Value = [double]"1.74439431788187E-08"
$ws.Cells.Item(25, 4).Value = 1
$ws.Cells.Item(26, 1).Value = 0.0002142307639587671
$ws.Cells.Item(26, 2).Value = 0.9999464750289917
$ws.Cells.Item(26, 3).Value = [double]"8.460473566174187E-08"
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(27, 1).Value = 0.0002401038655079901
$ws.Cells.Item(27, 2).Value = 0.9999464750289917
$ws.Cells.Item(27, 3).Value = [double]"1.976817465276781E-08"
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(28, 1).Value = [double]"2.996850162162445E-05"
$ws.Cells.Item(28, 2).Value = 1
$ws.Cells.Item(28, 3).Value = [double]"1.039930452861881E-08"
$ws.Cells.Item(28, 4).Value = 1
$ws.Cells.Item(29, 1).Value = 0.0002828243596013635
$ws.Cells.Item(29, 2).Value = 0.9999286532402039
$ws.Cells.Item(29, 3).Value = [double]"1.610216138203668E-08"
$ws.Cells.Item(29, 4).Value = 1
$ws.Cells.Item(30, 1).Value = 0.000713373941835016
$ws.Cells.Item(30, 2).Value = 0.9998751878738403
$ws.Cells.Item(30, 3).Value = [double]"1.967234020128217E-08"
$ws.Cells.Item(30, 4).Value = 1
$ws.Cells.Item(31, 1).Value = 0.00068421580363065
$ws.Cells.Item(31, 2).Value = 0.9998930096626282
$ws.Cells.Item(31, 3).Value = [double]"8.709644561122332E-08"
$ws.Cells.Item(31, 4).Value = 1
$ws.Cells.Item(32, 1).Value = 0.0001019663468468934
$ws.Cells.Item(32, 2).Value = 0.9999464750289917
$ws.Cells.Item(32, 3).Value = [double]"9.970057135433308E-08"
$ws.Cells.Item(32, 4).Value = 1
$ws.Cells.Item(33, 1).Value = 0.0004655737138818949
$ws.Cells.Item(33, 2).Value = 0.999910831451416
$ws.Cells.Item(33, 3).Value = [double]"8.132216322565E-08"
$ws.Cells.Item(33, 4).Value = 1
$ws.Cells.Item(34, 1).Value = 0.0004570386081468314
$ws.Cells.Item(34, 2).Value = 0.999910831451416
$ws.Cells.Item(34, 3).Value = [double]"1.141917209679377E-07"
$ws.Cells.Item(34, 4).Value = 1
$ws.Cells.Item(35, 1).Value = 0.0005027658771723509
$ws.Cells.Item(35, 2).Value = 0.9999464750289917
$ws.Cells.Item(35, 3).Value = [double]"4.013433141381029E-08"
$ws.Cells.Item(35, 4).Value = 1
$ws.Cells.Item(36, 1).Value = [double]"7.791847747284919E-05"
$ws.Cells.Item(36, 2).Value = 0.9999821782112122
$ws.Cells.Item(36, 3).Value = [double]"2.578201296898897E-08"
$ws.Cells.Item(36, 4).Value = 1
$ws.Cells.Item(37, 1).Value = 0.0001749493094393983
$ws.Cells.Item(37, 2).Value = 0.9999821782112122
$ws.Cells.Item(37, 3).Value = [double]"7.931264534022375E-09"
$ws.Cells.Item(37, 4).Value = 1
$ws.Cells.Item(38, 1).Value = [double]"5.001755562261678E-05"
$ws.Cells.Item(38, 2).Value = 0.9999821782112122
$ws.Cells.Item(38, 3).Value = [double]"6.900938043230553E-09"
$ws.Cells.Item(38, 4).Value = 1
$ws.Cells.Item(39, 1).Value = [double]"1.515305848442949E-05"
$ws.Cells.Item(39, 2).Value = 1
$ws.Cells.Item(39, 3).Value = [double]"3.666126957568849E-09"
$ws.Cells.Item(39, 4).Value = 1
$ws.Cells.Item(40, 1).Value = [double]"1.842087294789962E-05"
$ws.Cells.Item(40, 2).Value = 1
$ws.Cells.Item(40, 3).Value = [double]"8.146963104138649E-10"
$ws.Cells.Item(40, 4).Value = 1
$ws.Cells.Item(41, 1).Value = 0.0006934404955245554
$ws.Cells.Item(41, 2).Value = 0.9998930096626282
$ws.Cells.Item(41, 3).Value = [double]"1.408926397772348E-08"
$ws.Cells.Item(41, 4).Value = 1
$ws.Cells.Item(42, 1).Value = [double]"1.472023177484516E-05"
$ws.Cells.Item(42, 2).Value = 1
$ws.Cells.Item(42, 3).Value = [double]"6.445650235775702E-09"
$ws.Cells.Item(42, 4).Value = 1
$ws.Cells.Item(43, 1).Value = 0.0001574500056449324
$ws.Cells.Item(43, 2).Value = 0.9999643564224243
$ws.Cells.Item(43, 3).Value = [double]"3.043126417523467E-09"
$ws.Cells.Item(43, 4).Value = 1
$ws.Cells.Item(44, 1).Value = 0.0002356106997467577
$ws.Cells.Item(44, 2).Value = 0.9999464750289917
$ws.Cells.Item(44, 3).Value = [double]"1.293929186019227E-09"
$ws.Cells.Item(44, 4).Value = 1
$ws.Cells.Item(45, 1).Value = 0.0006225730176083744
$ws.Cells.Item(45, 2).Value = 0.999910831451416
$ws.Cells.Item(45, 3).Value = [double]"7.787487987798158E-09"
$ws.Cells.Item(45, 4).Value = 1
$ws.Cells.Item(46, 1).Value = 0.0002640956663526595
$ws.Cells.Item(46, 2).Value = 0.9999286532402039
$ws.Cells.Item(46, 3).Value = [double]"2.707661650447335E-09"
$ws.Cells.Item(46, 4).Value = 1
$ws.Cells.Item(47, 1).Value = 0.0003335273650009185
$ws.Cells.Item(47, 2).Value = 0.999910831451416
$ws.Cells.Item(47, 3).Value = [double]"8.098941961520723E-09"
$ws.Cells.Item(47, 4).Value = 1
$ws.Cells.Item(48, 1).Value = [double]"1.525215884612408E-05"
$ws.Cells.Item(48, 2).Value = 1
$ws.Cells.Item(48, 3).Value = [double]"3.306691365168035E-09"
$ws.Cells.Item(48, 4).Value = 1
$ws.Cells.Item(49, 1).Value = 0.0003151391283608973
$ws.Cells.Item(49, 2).Value = 0.999910831451416
$ws.Cells.Item(49, 3).Value = [double]"3.718592367363271E-08"
$ws.Cells.Item(49, 4).Value = 1
$ws.Cells.Item(50, 1).Value = [double]"2.849831435014494E-05"
$ws.Cells.Item(50, 2).Value = 0.9999821782112122
$ws.Cells.Item(50, 3).Value = [double]"6.277922182107432E-09"
$ws.Cells.Item(50, 4).Value = 1
$ws.Cells.Item(51, 1).Value = [double]"2.153246896341443E-05"
$ws.Cells.Item(51, 2).Value = 1
$ws.Cells.Item(51, 3).Value = [double]"3.809873305726796E-09"
$ws.Cells.Item(51, 4).Value = 1
